# Weekly update: insert a new market-report row (price observation) for
# "Pepino ensalada" (Vega Monumental Concepción) dated 2021-09-30, pushing
# the existing rows 21..72 down to 22..73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (shifts rows 21:72 down to 22:73,
# carrying formatting - e.g. the date style - from the row above).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new observation.
$ws.Cells.Item(21, 1).Value2  = 11
$ws.Cells.Item(21, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value2  = "Bíobío"
$ws.Cells.Item(21, 4).Value2  = 44469
$ws.Cells.Item(21, 5).Value2  = 8
$ws.Cells.Item(21, 6).Value2  = 100112043
$ws.Cells.Item(21, 7).Value2  = "Pepino ensalada"
$ws.Cells.Item(21, 8).Value2  = "Sin especificar"
$ws.Cells.Item(21, 9).Value2  = "Primera"
$ws.Cells.Item(21, 10).Value2 = 100
$ws.Cells.Item(21, 11).Value2 = 16000
$ws.Cells.Item(21, 12).Value2 = 17000
$ws.Cells.Item(21, 13).Value2 = 16500
$ws.Cells.Item(21, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(21, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value2 = 275
$ws.Cells.Item(21, 17).Value2 = 60
$ws.Cells.Item(21, 18).Value2 = "Hortaliza"
